# Apply cryptos list update (prices / volume% refresh) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.227.87"
$ws.Range("E2").Value = "  -4.67%  "

$ws.Range("D3").Value = "3.087.71"
$ws.Range("E3").Value = "  -4.51%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'549.83"
$ws.Range("E5").Value = "  -4.41%  "

$ws.Range("D6").Value = "'136.18"
$ws.Range("E6").Value = "  -11.89%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "3.081.28"
$ws.Range("E8").Value = "  -4.45%  "

$ws.Range("E9").Value = "  -3.45%  "

$ws.Range("E10").Value = "  -5.92%  "

$ws.Range("D11").Value = "'6.33"
$ws.Range("E11").Value = "  -10.38%  "

$ws.Range("D12").Value = "'0.466"
$ws.Range("E12").Value = "  -3.61%  "

$ws.Range("D13").Value = "'35.24"
$ws.Range("E13").Value = "  -7.06%  "

$ws.Range("E14").Value = "  -7.23%  "

$ws.Range("D15").Value = "3.585.76"
$ws.Range("E15").Value = "  -4.44%  "

$ws.Range("D16").Value = "63.303.32"
$ws.Range("E16").Value = "  -4.57%  "

$ws.Range("E17").Value = "  -3.12%  "

$ws.Range("D18").Value = "3.082.52"
$ws.Range("E18").Value = "  -4.59%  "

$ws.Range("D19").Value = "'6.71"
$ws.Range("E19").Value = "  -5.08%  "

$ws.Range("D20").Value = "'490.97"
$ws.Range("E20").Value = "  -11.53%  "

$ws.Range("D21").Value = "'13.61"
$ws.Range("E21").Value = "  -5.12%  "

$ws.Range("D22").Value = "'0.712"
$ws.Range("E22").Value = "  -3.47%  "

$ws.Range("D23").Value = "'7.23"
$ws.Range("E23").Value = "  -7.76%  "

$ws.Range("D24").Value = "'78.49"
$ws.Range("E24").Value = "  -4.13%  "

$ws.Range("D25").Value = "'12.35"
$ws.Range("E25").Value = "  -8.78%  "

$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("D27").Value = "'8.47"
$ws.Range("E27").Value = "  -10.18%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("D30").Value = "'1.98"
$ws.Range("E30").Value = "  -12.67%  "

$ws.Range("D31").Value = "'26.52"
$ws.Range("E31").Value = "  -4.38%  "

$ws.Range("D32").Value = "'1.13"
$ws.Range("E32").Value = "  -3.70%  "

$ws.Range("D33").Value = "'2.51"
$ws.Range("E33").Value = "  -8.74%  "

$ws.Range("D34").Value = "'58.73"

$ws.Range("D35").Value = "'522.37"
$ws.Range("E35").Value = "  -7.17%  "

$ws.Range("E36").Value = "  -6.31%  "

$ws.Range("D37").Value = "'5.13"
$ws.Range("E37").Value = "  -10.65%  "

$ws.Range("D38").Value = "'0.0405"
$ws.Range("E38").Value = "  -12.37%  "

$ws.Range("D39").Value = "3.130.29"
$ws.Range("E39").Value = "  +0.08%  "

$ws.Range("E40").Value = "  -7.81%  "

$ws.Range("D41").Value = "'0.119"
$ws.Range("E41").Value = "  -5.63%  "

$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").Value = "'8.14"
$ws.Range("E42").Value = "  -5.51%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.67"
$ws.Range("E43").Value = "  -11.50%  "

$ws.Range("E44").Value = "  -6.04%  "

$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").Value = "'2.08"
$ws.Range("E46").Value = "  -10.87%  "

$ws.Range("D47").Value = "'24.89"
$ws.Range("E47").Value = "  -7.96%  "

$ws.Range("D48").Value = "'121.54"
$ws.Range("E48").Value = "  -0.98%  "

$ws.Range("E49").Value = "  -4.21%  "

$ws.Range("E50").Value = "  -10.43%  "

$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'2.03"
$ws.Range("E51").Value = "  -9.48%  "
